# Update "想去人数" (number of people interested) values in column F
# for both the "展览" and "全部类型" worksheets, per the source diff.

$wb = $excel.ActiveWorkbook

# Row -> new F value mapping (applies identically to both data sheets)
$updates = @{
    2  = 2022
    7  = 1684
    8  = 25
    13 = 99
    17 = 116
    19 = 3903
    22 = 440
    23 = 365
    24 = 826
    25 = 551
    28 = 1697
    31 = 167
    32 = 13
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
